$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins, Losses, Ties) in AD1:AF1,
# matching the style of the existing header row (e.g. A1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill season-record columns for every data row (2-54) with the
# team's win/loss/tie totals for the season.
for ($row = 2; $row -le 54; $row++) {
    $ws.Cells.Item($row, 30).Value = 95
    $ws.Cells.Item($row, 31).Value = 67
    $ws.Cells.Item($row, 32).Value = 0
}
